# Apply the "Add files via upload" revision to ImageConversionTestCases.xlsx
# (test case renumbering / rewording + a couple of row-height and
# view/selection tweaks on Sheet1).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Text edits (only cells whose displayed text actually changed) ---

$ws.Range("C8").Value  = "Load file with invald image file format"
$ws.Range("C19").Value = "Load file with invald image file format"
$ws.Range("C23").Value = "Load file with invald image file format"

$ws.Range("C11").Value = "Test invalid image width with valid height"
$ws.Range("C14").Value = "Test invalid height with valid width"
$ws.Range("C15").Value = "Test very large height with valid width"

$ws.Range("C27").Value = "Input non empty matrix filled with 0s"
$ws.Range("C28").Value = "Input non empty matrix filled with 1s"

$ws.Range("C29").Value = "Input non empty matrix with  number of rows > 49 filled with 0s"
$ws.Range("D29").Value = "matrix = 49 rows filled with 0s"

$ws.Range("C30").Value = "Input non empty matrix with number of columns > 251 filled with 1s"

$ws.Range("C31").Value = "Input non empty matrix filled with integers other than 0 or 1"
$ws.Range("C32").Value = "Input non empty matrix filled with not integers"

$ws.Range("C35").Value = "Input nominal width"
$ws.Range("C36").Value = "Input width >= 250"

# --- Row height tweaks (wrapped Test-Steps column needs more room) ---

$ws.Rows.Item(11).RowHeight = 30
$ws.Rows.Item(14).RowHeight = 30
$ws.Rows.Item(15).RowHeight = 30
$ws.Rows.Item(27).RowHeight = 30
$ws.Rows.Item(28).RowHeight = 30
$ws.Rows.Item(29).RowHeight = 45
$ws.Rows.Item(30).RowHeight = 45

# --- View/selection state ---

$ws.Range("D36").Select()
